# Applies market-data value refresh produced by the scheduled Moogle Profits runner.
# Updates currentAveragePrice* / Leve price / profit columns (H-N) across all job sheets.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4681.148
$ws.Range("I64").Value = 4995.75
$ws.Range("J64").Value = 4223.5454
$ws.Range("K64").Value = 4995.75
$ws.Range("L64").Value = 4223.5454
$ws.Range("M64").Value = -4747.75
$ws.Range("N64").Value = -4719.5454
$ws.Range("H67").Value = 4681.148
$ws.Range("I67").Value = 4995.75
$ws.Range("J67").Value = 4223.5454
$ws.Range("K67").Value = 4995.75
$ws.Range("L67").Value = 4223.5454
$ws.Range("M67").Value = -4137.75
$ws.Range("N67").Value = -5939.5454
$ws.Range("H100").Value = 2399.3684
$ws.Range("I100").Value = 2293.5881
$ws.Range("K100").Value = 2293.5881
$ws.Range("M100").Value = -1752.5881
$ws.Range("H137").Value = 3591.2856
$ws.Range("I137").Value = 2457.4
$ws.Range("J137").Value = 4221.222
$ws.Range("K137").Value = 7372.200000000001
$ws.Range("L137").Value = 12663.666
$ws.Range("M137").Value = -4822.200000000001
$ws.Range("N137").Value = -17763.666

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 834.2
$ws.Range("J4").Value = 891
$ws.Range("L4").Value = 891
$ws.Range("N4").Value = -1123
$ws.Range("H32").Value = 14963.725
$ws.Range("I32").Value = 6734.8184
$ws.Range("K32").Value = 6734.8184
$ws.Range("M32").Value = -6447.8184
$ws.Range("H61").Value = 26323132
$ws.Range("I61").Value = 5899.933
$ws.Range("J61").Value = 125012750
$ws.Range("K61").Value = 5899.933
$ws.Range("L61").Value = 125012750
$ws.Range("M61").Value = -5687.933
$ws.Range("N61").Value = -125013174
$ws.Range("H74").Value = 5400.0625
$ws.Range("I74").Value = 2114
$ws.Range("J74").Value = 8686.125
$ws.Range("K74").Value = 2114
$ws.Range("L74").Value = 8686.125
$ws.Range("M74").Value = -1240
$ws.Range("N74").Value = -10434.125
$ws.Range("H77").Value = 5400.0625
$ws.Range("I77").Value = 2114
$ws.Range("J77").Value = 8686.125
$ws.Range("K77").Value = 10570
$ws.Range("L77").Value = 43430.625
$ws.Range("M77").Value = -6202
$ws.Range("N77").Value = -52166.625
$ws.Range("H88").Value = 1792.875
$ws.Range("I88").Value = 1085.8572
$ws.Range("J88").Value = 2342.7778
$ws.Range("K88").Value = 1085.8572
$ws.Range("L88").Value = 2342.7778
$ws.Range("M88").Value = -679.8571999999999
$ws.Range("N88").Value = -3154.7778
$ws.Range("H91").Value = 1792.875
$ws.Range("I91").Value = 1085.8572
$ws.Range("J91").Value = 2342.7778
$ws.Range("K91").Value = 1085.8572
$ws.Range("L91").Value = 2342.7778
$ws.Range("M91").Value = 318.1428000000001
$ws.Range("N91").Value = -5150.7778
$ws.Range("H102").Value = 1517.4546
$ws.Range("I102").Value = 1361
$ws.Range("K102").Value = 1361
$ws.Range("M102").Value = 261
$ws.Range("H136").Value = 26323132
$ws.Range("I136").Value = 5899.933
$ws.Range("J136").Value = 125012750
$ws.Range("K136").Value = 17699.799
$ws.Range("L136").Value = 375038250
$ws.Range("M136").Value = -15149.799
$ws.Range("N136").Value = -375043350

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 50004.5
$ws.Range("I19").Value = 5009
$ws.Range("K19").Value = 5009
$ws.Range("M19").Value = -4836
$ws.Range("H99").Value = 2051.2307
$ws.Range("I99").Value = 1805.5
$ws.Range("K99").Value = 1805.5
$ws.Range("M99").Value = -307.5
$ws.Range("H134").Value = 3857.1428
$ws.Range("I134").Value = 1500
$ws.Range("K134").Value = 4500
$ws.Range("M134").Value = -1965

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 32564.889
$ws.Range("I36").Value = 9198.333000000001
$ws.Range("J36").Value = 44248.168
$ws.Range("K36").Value = 9198.333000000001
$ws.Range("L36").Value = 44248.168
$ws.Range("M36").Value = -8810.333000000001
$ws.Range("N36").Value = -45024.168
$ws.Range("H38").Value = 20597.8
$ws.Range("I38").Value = 5500
$ws.Range("J38").Value = 30663
$ws.Range("K38").Value = 5500
$ws.Range("L38").Value = 30663
$ws.Range("M38").Value = -5123
$ws.Range("N38").Value = -31417
$ws.Range("H39").Value = 26340.25
$ws.Range("I39").Value = 19170.715
$ws.Range("J39").Value = 36377.6
$ws.Range("K39").Value = 19170.715
$ws.Range("L39").Value = 36377.6
$ws.Range("M39").Value = -18779.715
$ws.Range("N39").Value = -37159.6
$ws.Range("H40").Value = 32564.889
$ws.Range("I40").Value = 9198.333000000001
$ws.Range("J40").Value = 44248.168
$ws.Range("K40").Value = 9198.333000000001
$ws.Range("L40").Value = 44248.168
$ws.Range("M40").Value = -9038.333000000001
$ws.Range("N40").Value = -44568.168
$ws.Range("H44").Value = 50000
$ws.Range("J44").Value = 50000
$ws.Range("L44").Value = 50000
$ws.Range("N44").Value = -50884
$ws.Range("H46").Value = 20597.8
$ws.Range("I46").Value = 5500
$ws.Range("J46").Value = 30663
$ws.Range("K46").Value = 5500
$ws.Range("L46").Value = 30663
$ws.Range("M46").Value = -5289
$ws.Range("N46").Value = -31085
$ws.Range("H47").Value = 34999.5
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("H49").Value = 26340.25
$ws.Range("I49").Value = 19170.715
$ws.Range("J49").Value = 36377.6
$ws.Range("K49").Value = 19170.715
$ws.Range("L49").Value = 36377.6
$ws.Range("M49").Value = -18988.715
$ws.Range("N49").Value = -36741.6
$ws.Range("H58").Value = 3436.8262
$ws.Range("I58").Value = 2026.65
$ws.Range("K58").Value = 2026.65
$ws.Range("M58").Value = -1823.65
$ws.Range("H136").Value = 3436.8262
$ws.Range("I136").Value = 2026.65
$ws.Range("K136").Value = 6079.950000000001
$ws.Range("M136").Value = -3529.950000000001

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 927.3333
$ws.Range("I5").Value = 600
$ws.Range("J5").Value = 992.8
$ws.Range("K5").Value = 1800
$ws.Range("L5").Value = 2978.4
$ws.Range("M5").Value = -1688
$ws.Range("N5").Value = -3202.4
$ws.Range("H135").Value = 927.3333
$ws.Range("I135").Value = 600
$ws.Range("J135").Value = 992.8
$ws.Range("K135").Value = 5400
$ws.Range("L135").Value = 8935.199999999999
$ws.Range("M135").Value = -2865
$ws.Range("N135").Value = -14005.2
$ws.Range("H140").Value = 1536.7587
$ws.Range("J140").Value = 1650.5
$ws.Range("L140").Value = 4951.5
$ws.Range("N140").Value = -15311.5

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2771.5
$ws.Range("I113").Value = 1791.6666
$ws.Range("K113").Value = 1791.6666
$ws.Range("M113").Value = 378.3334

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7635.7827
$ws.Range("I132").Value = 7635.7827
$ws.Range("K132").Value = 22907.3481
$ws.Range("M132").Value = -20377.3481
$ws.Range("H136").Value = 16672772
$ws.Range("I136").Value = 3822.5
$ws.Range("K136").Value = 11467.5
$ws.Range("M136").Value = -8917.5

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 8504
$ws.Range("J45").Value = 9631.25
$ws.Range("L45").Value = 9631.25
$ws.Range("N45").Value = -10613.25
$ws.Range("H126").Value = 2598
$ws.Range("I126").Value = 2442.9092
$ws.Range("K126").Value = 7328.7276
$ws.Range("M126").Value = -4858.7276
$ws.Range("H132").Value = 5664
$ws.Range("I132").Value = 1999.6666
$ws.Range("J132").Value = 9328.333000000001
$ws.Range("K132").Value = 5998.9998
$ws.Range("L132").Value = 27984.999
$ws.Range("M132").Value = -3468.9998
$ws.Range("N132").Value = -33044.999

Write-Output "Moogle_Profits value refresh applied"